$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J19").Value = 12.434613462352335
$ws.Range("J20").Value = 16.80050595536094
$ws.Range("J21").Value = 11.282963378125267
$ws.Range("J22").Value = 25.042808754677555
$ws.Range("J23").Value = 3.2011163356916352
$ws.Range("J24").Value = 13.523574517571838
$ws.Range("J25").Value = 6.1196997869329204
$ws.Range("J26").Value = 5.9488136666578013
$ws.Range("J27").Value = 5.2451982064110645

$ws.Range("N8").Select()
